# Fixed #253 Moving from POI 3.16 to 3.17.
#
# The document embeds a captured Java stack trace (as literal text,
# inside a single run) illustrating a failing asTable(...) call. Bumping
# the POI dependency shifted the M2Doc/Acceleo call-site line numbers
# reported in that trace, and one pair of now-inlined/removed stack
# frames (the old caseTemplate(...) frames feeding into
# caseDocumentTemplate(...)) disappeared entirely from the trace.
#
# Rebuild the "old" and "new" flavors of the affected stack-trace
# excerpt (one array entry per source line) and swap them with a single
# Find/Replace over the whole block so every frame line lands exactly
# where the new trace puts it.

$oldLines = @(
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "	at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)",
    "	at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:516)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)",
    "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:172)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1158)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)",
    "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:183)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:311)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:1)",
    "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:201)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:266)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)",
    "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:246)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:255)",
    "	at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:705)",
    "	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:458)",
    "	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:365)",
    "	at sun.reflect.GeneratedMethodAccessor76.invoke(Unknown Source)",
    "	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)",
    "	at java.lang.reflect.Method.invoke(Method.java:498)",
    "	at org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:50)"
)

$newLines = @(
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "	at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)",
    "	at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:540)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)",
    "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)",
    "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)",
    "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)",
    "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)",
    "	at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)",
    "	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)",
    "	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)",
    "	at sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)",
    "	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)",
    "	at java.lang.reflect.Method.invoke(Method.java:498)",
    "	at org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:50)"
)

$oldStack = $oldLines -join "`n"
$newStack = $newLines -join "`n"

$d = $word.ActiveDocument
$found = $d.Content.Find.Execute($oldStack, $false, $false, $false, $false, $false, $true, 1, $false, $newStack, 2)
if (-not $found) {
    throw "Could not locate the stack trace block to update"
}
